$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 127
$ws1.Range("F4").Value = 2062
$ws1.Range("F5").Value = 350
$ws1.Range("F6").Value = 599
$ws1.Range("F9").Value = 10589
$ws1.Range("F10").Value = 179
$ws1.Range("F11").Value = 156
$ws1.Range("F14").Value = 413
$ws1.Range("F15").Value = 7475
$ws1.Range("F18").Value = 230
$ws1.Range("F20").Value = 3318

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 127
$ws4.Range("F4").Value = 2062
$ws4.Range("F5").Value = 350
$ws4.Range("F6").Value = 599
$ws4.Range("F12").Value = 10589
$ws4.Range("F13").Value = 179
$ws4.Range("F14").Value = 156
$ws4.Range("F17").Value = 413
$ws4.Range("F18").Value = 7475
$ws4.Range("F21").Value = 230
$ws4.Range("F23").Value = 3318
